$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 14 (pushing existing rows 14+ down to 16+),
# to add two new "Maternal Characteristics" covariates: teen_birth and geriatric_birth.
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(14).Insert()

# Column A (Name) for the two new rows
$ws.Cells.Item(14, 1).Value2 = "teen_birth"
$ws.Cells.Item(15, 1).Value2 = "geriatric_birth"

# Column B (Type)
$ws.Cells.Item(14, 2).Value2 = "Covariate"
$ws.Cells.Item(15, 2).Value2 = "Covariate"

# Column C (Category)
$ws.Cells.Item(14, 3).Value2 = "Maternal Characteristics"
$ws.Cells.Item(15, 3).Value2 = "Maternal Characteristics"

# Column D (Source)
$ws.Cells.Item(14, 4).Value2 = "Census Reporter"
$ws.Cells.Item(15, 4).Value2 = "Census Reporter"

# Column E (Description)
$ws.Cells.Item(14, 5).Value2 = "proportion of births with maternal age 15-19"
$ws.Cells.Item(15, 5).Value2 = "proportion of births with maternal age 35+"

# Restore the selected cell as it ended up after the edit
$ws.Range("E16").Select() | Out-Null
